$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3251.1667
$ws.Range("J40").Value = 5002
$ws.Range("L40").Value = 5002
$ws.Range("N40").Value = -5352
$ws.Range("H51").Value = 34420.31
$ws.Range("I51").Value = 7500
$ws.Range("K51").Value = 7500
$ws.Range("M51").Value = -7016
$ws.Range("H62").Value = 44450930
$ws.Range("I62").Value = 53339910
$ws.Range("J62").Value = 6000
$ws.Range("K62").Value = 53339910
$ws.Range("L62").Value = 6000
$ws.Range("M62").Value = -53339286
$ws.Range("N62").Value = -7248
$ws.Range("H65").Value = 44450930
$ws.Range("I65").Value = 53339910
$ws.Range("J65").Value = 6000
$ws.Range("K65").Value = 266699550
$ws.Range("L65").Value = 30000
$ws.Range("M65").Value = -266696430
$ws.Range("N65").Value = -36240
$ws.Range("H92").Value = 1223.0322
$ws.Range("I92").Value = 1308.5416
$ws.Range("K92").Value = 1308.5416
$ws.Range("M92").Value = -60.54160000000002
$ws.Range("H137").Value = 3161.7273
$ws.Range("I137").Value = 2563.4666
$ws.Range("K137").Value = 7690.399800000001
$ws.Range("M137").Value = -5140.399800000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("H61").Value = 6369.1665
$ws.Range("I61").Value = 4770.231
$ws.Range("K61").Value = 4770.231
$ws.Range("M61").Value = -4558.231
$ws.Range("H110").Value = 11906380
$ws.Range("I110").Value = 16667532
$ws.Range("J110").Value = 3500
$ws.Range("K110").Value = 16667532
$ws.Range("L110").Value = 3500
$ws.Range("M110").Value = -16665487
$ws.Range("N110").Value = -7590
$ws.Range("H125").Value = 23999
$ws.Range("J125").Value = 23999
$ws.Range("L125").Value = 23999
$ws.Range("N125").Value = -33839
$ws.Range("H136").Value = 6369.1665
$ws.Range("I136").Value = 4770.231
$ws.Range("K136").Value = 14310.693
$ws.Range("M136").Value = -11760.693

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H34").Value = 5250
$ws.Range("J34").Value = 9500
$ws.Range("L34").Value = 9500
$ws.Range("N34").Value = -9728
$ws.Range("H37").Value = 1200
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()
$ws.Range("H99").Value = 1691.2084
$ws.Range("J99").Value = 1850
$ws.Range("L99").Value = 1850
$ws.Range("N99").Value = -4846
$ws.Range("H105").Value = 38477080
$ws.Range("I105").Value = 71452136
$ws.Range("J105").Value = 6186.4165
$ws.Range("K105").Value = 71452136
$ws.Range("L105").Value = 6186.4165
$ws.Range("M105").Value = -71450389
$ws.Range("N105").Value = -9680.416499999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H55").Value = 29693.8
$ws.Range("I55").Value = 26485.5
$ws.Range("J55").Value = 31832.666
$ws.Range("K55").Value = 26485.5
$ws.Range("L55").Value = 31832.666
$ws.Range("M55").Value = -26170.5
$ws.Range("N55").Value = -32462.666
$ws.Range("H134").Value = 4831.68
$ws.Range("J134").Value = 6983.1
$ws.Range("L134").Value = 20949.3
$ws.Range("N134").Value = -26019.3

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 57179950
$ws.Range("I4").Value = 67780470
$ws.Range("K4").Value = 203341410
$ws.Range("M4").Value = -203341298
$ws.Range("H37").Value = 87562.08
$ws.Range("J37").Value = 87562.08
$ws.Range("L37").Value = 262686.24
$ws.Range("N37").Value = -262910.24
$ws.Range("H56").Value = 11375.5
$ws.Range("I56").Value = 11375.5
$ws.Range("K56").Value = 11375.5
$ws.Range("M56").Value = -10845.5
$ws.Range("H58").Value = 2508
$ws.Range("I58").Value = 2622.5
$ws.Range("J58").Value = 2493.6875
$ws.Range("K58").Value = 7867.5
$ws.Range("L58").Value = 7481.0625
$ws.Range("M58").Value = -7739.5
$ws.Range("N58").Value = -7737.0625
$ws.Range("H98").Value = 350
$ws.Range("I98").Value = 200
$ws.Range("J98").Value = 500
$ws.Range("K98").Value = 600
$ws.Range("L98").Value = 1500
$ws.Range("M98").Value = 898
$ws.Range("N98").Value = -4496

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 49961.668
$ws.Range("J93").Value = 49961.668
$ws.Range("L93").Value = 49961.668
$ws.Range("N93").Value = -53705.668
$ws.Range("H102").Value = 992.5
$ws.Range("I102").Value = 992.5
$ws.Range("K102").Value = 992.5
$ws.Range("M102").Value = 629.5
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H126").Value = 3207.6667
$ws.Range("I126").Value = 2440
$ws.Range("K126").Value = 7320
$ws.Range("M126").Value = -4850

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3187.375
$ws.Range("J7").Value = 5969.6665
$ws.Range("L7").Value = 5969.6665
$ws.Range("N7").Value = -6193.6665
$ws.Range("H46").Value = 6435.6
$ws.Range("I46").Value = 1754.2
$ws.Range("J46").Value = 7371.88
$ws.Range("K46").Value = 1754.2
$ws.Range("L46").Value = 7371.88
$ws.Range("M46").Value = -1566.2
$ws.Range("N46").Value = -7747.88
$ws.Range("H100").Value = 8067847
$ws.Range("I100").Value = 35716110
$ws.Range("K100").Value = 35716110
$ws.Range("M100").Value = -35715569
$ws.Range("H126").Value = 3187.375
$ws.Range("J126").Value = 5969.6665
$ws.Range("L126").Value = 17908.9995
$ws.Range("N126").Value = -22848.9995
$ws.Range("H136").Value = 5260.0415
$ws.Range("J136").Value = 5982.25
$ws.Range("L136").Value = 17946.75
$ws.Range("N136").Value = -23046.75
$ws.Range("H140").Value = 100000
$ws.Range("J140").Value = 100000
$ws.Range("L140").Value = 100000
$ws.Range("N140").Value = -110360

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 5075
$ws.Range("H31").Value = 80000
$ws.Range("I31").Value = 80000
$ws.Range("K31").Value = 80000
$ws.Range("M31").Value = -79652
$ws.Range("H76").Value = 61999.668
$ws.Range("J76").Value = 61999.668
$ws.Range("L76").Value = 61999.668
$ws.Range("N76").Value = -62629.668
$ws.Range("H79").Value = 61999.668
$ws.Range("J79").Value = 61999.668
$ws.Range("L79").Value = 61999.668
$ws.Range("N79").Value = -64183.668
$ws.Range("H96").Value = 1478
$ws.Range("J96").Value = 1529.6
$ws.Range("L96").Value = 1529.6
$ws.Range("N96").Value = -4275.6
$ws.Range("H107").Value = 4359.6
$ws.Range("I107").Value = 1119.8
$ws.Range("J107").Value = 7599.4
$ws.Range("K107").Value = 3359.4
$ws.Range("L107").Value = 22798.2
$ws.Range("M107").Value = -1439.4
$ws.Range("N107").Value = -26638.2
$ws.Range("H126").Value = 5128.773
$ws.Range("I126").Value = 4956.6
$ws.Range("J126").Value = 5497.7144
$ws.Range("K126").Value = 14869.8
$ws.Range("L126").Value = 16493.1432
$ws.Range("M126").Value = -12399.8
$ws.Range("N126").Value = -21433.1432
$ws.Range("H136").Value = 3516.6296
$ws.Range("I136").Value = 1391.8572
$ws.Range("J136").Value = 5804.846
$ws.Range("K136").Value = 4175.571599999999
$ws.Range("L136").Value = 17414.538
$ws.Range("M136").Value = -1625.571599999999
$ws.Range("N136").Value = -22514.538
